{"js": "// Add \"Crossplane\" to the DevOps skills line, right after \"Ansible, Terraform, \"\n// and before \"DataDog\" -- i.e. \"... Ansible, Terraform, Crossplane, DataDog, ...\"\nconst results = context.document.body.search(\"Ansible, Terraform, \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Ansible, Terraform, \" in the document body.');\n}\n\n// Insert the new skill (with its trailing \", \") immediately after the matched\n// text, right before \"DataDog\", inheriting formatting from the adjacent run.\nresults.items[0].insertText(\"Crossplane, \", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add \"Crossplane\" to the DevOps skills line, right after \"Ansible, Terraform, \"\n# and before \"DataDog\" -- i.e. \"... Ansible, Terraform, Crossplane, DataDog, ...\"\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Ansible, Terraform, \"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"Ansible, Terraform, \" in the document.'\n}\n\n# Collapse to the end of the match so the insertion lands right before \"DataDog\",\n# inheriting formatting from the text it is adjacent to.\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\"Crossplane, \")\n"}
